$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: swap E11/F11
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = "['MEC-2A-Ajustagem', -, -, -]"

# Row 12: swap E12/F12
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "['MEC-2A-Ajustagem', -, -, -]"

# Row 14: swap E14/F14
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "['MEC-2A-Ajustagem', -, -, -]"

# Row 15: swap E15/F15
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "['MEC-2A-Ajustagem', -, -, -]"

# Row 18
$ws.Range("B18").Value = "-"
$ws.Range("D18").Value = "[-, -, -, 'ELM-1NA-Processos de Usinagem 1']"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "[-, 'MEC-1NA-Tornearia', -, 'MEC-1NA-Metrologia 1']"

# Row 19
$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "[-, -, -, 'MEC-1NB-Tornearia']"
$ws.Range("D19").Value = "[-, -, -, 'ELM-1NA-Processos de Usinagem 1']"
$ws.Range("E19").Value = "-"
$ws.Range("F19").Value = "[-, 'MEC-1NA-Tornearia', -, 'MEC-1NA-Metrologia 1']"

# Row 20
$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "[-, -, 'MEC-1NB-Tornearia', -]"
$ws.Range("D20").Value = "[-, 'MEC-1NB-Tornearia', -, 'ELM-1NA-Processos de Usinagem 1']"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "[-, 'MEC-1NA-Tornearia', -, 'MEC-1NA-Metrologia 1']"

# Row 21
$ws.Range("B21").Value = "[-, -, 'MEC-1NB-Tornearia', -]"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "[-, -, -, 'ELM-1NA-Processos de Usinagem 1']"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "[-, 'MEC-1NA-Tornearia', -, 'MEC-1NA-Metrologia 1']"
